$wb = $excel.ActiveWorkbook
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 368324
$excel.CalculateFullRebuild()
